# Append three new paragraphs (styled "Paragrafoelenco") at the very end of the
# document body, right after the existing last paragraph
# ("All'inizializzazione della prima immagine...") and before the section
# properties (sectPr):
#   1. an empty paragraph (paragraph-mark formatting only, no run)
#   2. a paragraph containing "30.11:"
#   3. a paragraph containing the long note about fixing the StackOverFlow bug
$d = $word.ActiveDocument

# Collapse a range to the very end of the document's main story content, i.e.
# right after the last paragraph mark and before the sectPr.
$insertPoint = $d.Range($d.Content.End, $d.Content.End)

$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:pStyle w:val="Paragrafoelenco"/>
<w:rPr>
<w:rFonts w:cstheme="minorHAnsi"/>
</w:rPr>
</w:pPr>
</w:p>
<w:p>
<w:pPr>
<w:pStyle w:val="Paragrafoelenco"/>
<w:rPr>
<w:rFonts w:cstheme="minorHAnsi"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:rFonts w:cstheme="minorHAnsi"/>
</w:rPr>
<w:t>30.11:</w:t>
</w:r>
</w:p>
<w:p>
<w:pPr>
<w:pStyle w:val="Paragrafoelenco"/>
<w:rPr>
<w:rFonts w:cstheme="minorHAnsi"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:rFonts w:cstheme="minorHAnsi"/>
</w:rPr>
<w:t>&#8220;Fixed&#8221; il System.StackOverFlow dovuto alla ricerca del valore 1631614771990 inserendo un counter, un controllo fintanto che non arriva al 100, quando avviene l&#8217;offset viene moltiplicato di 7.</w:t>
</w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$insertPoint.InsertXML($xml)
